# CCO_eCoaching_Log_Runbook.docx edit script
# Implements:
#   1) "25205 - New Submission: ..." TFS item rewritten to the new
#      "25431 - Historical page: add searching by reason/sub reason;" item.
#   2) A new log row added to the revision-history table documenting the
#      10/13/2022 "TFS 25431 - Historical page: ..." change by Lili Huang.
#   3) The bold changeset number near "Get the following from TFS ..."
#      updated from 51867 to 52098.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Replace the old "TFS 25205 - New Submission..." paragraph text with
#    the new "25431 - Historical page..." text (single run).
# ---------------------------------------------------------------------
$old1 = "25205 – New Submission: display log name with success message for single log submission."
$new1 = "25431 - Historical page: add searching by reason/sub reason;"
$found1 = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) {
    throw "Could not find the old TFS 25205 paragraph text to replace."
}

# ---------------------------------------------------------------------
# 2) Add a new row to the change-log table (2nd table in the document)
#    for the 10/13/2022 TFS 25431 entry, right after the existing
#    08/24/2022 row (it becomes the new last row of the table).
# ---------------------------------------------------------------------
$table = $d.Tables.Item(2)
$newRow = $table.Rows.Add()
$rowIdx = $table.Rows.Count
$row = $table.Rows.Item($rowIdx)

$row.Cells.Item(1).Range.Text = "10/13/2022"
$row.Cells.Item(2).Range.Text = "TFS 25431 - Historical page: add searching by reason/sub reason;"
$row.Cells.Item(3).Range.Text = "Lili Huang"

# ---------------------------------------------------------------------
# 3) Update the bold "Changeset" number from 51867 to 52098.
# ---------------------------------------------------------------------
$found2 = $d.Content.Find.Execute("51867", $false, $false, $false, $false, $false, $true, 1, $false, "52098", 2)
if (-not $found2) {
    throw "Could not find the Changeset 51867 number to replace."
}
